$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Avverkningsanmälningar")

# Column C ("Förändrad") rows 2 through 32 need to go from 45601 to 45602
for ($r = 2; $r -le 32; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45601) {
        $cell.Value = 45602
    }
}
